# Adds hydrogen (Pumped Hydro Storage, EN_STG_PHS) to the ELC_Storage sheet:
# - Row 6: new technology-set summary row for EN_STG_PHS
# - Rows 17-20 (the Commodities/Csets table) shift down to rows 23-26 to make
#   room for the new technology's detail block
# - Rows 14-16: new detail rows for EN_STG_PHS (mirrors the existing
#   EN_STG_4hBatt block in rows 10-12)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELC_Storage")

# Reusable custom number formats already present in the workbook so that we
# land on the same *visual* formatting as the existing rows use (Excel will
# reuse the matching style record rather than minting a new one).
$TEXT_FMT = "\Te\x\t"
$DEC_FMT = "0.00"

# --- Make room: push the Csets table (rows 17-20) down to rows 23-26 -------
$ws.Range("A17:A22").EntireRow.Insert()

# --- Row 6: register the new EN_STG_PHS technology (mirrors rows 4 & 5) ---
$ws.Range("B6").Value = "STGTSS"
$ws.Range("C6").Value = "EN_STG_PHS"
$ws.Range("D6").Value = "Pumped hydro storage"
$ws.Range("E6").Value = "TBTU"
$ws.Range("F6").Value = "GW"
$ws.Range("G6").Value = "DAYNITE"
$ws.Range("H6").NumberFormat = $TEXT_FMT

# --- Row 14: EN_STG_PHS main parameter row (mirrors row 10) ---------------
$ws.Range("B14").Formula = "=C6"
$ws.Range("C14").Value = "ELC"
$ws.Range("E14").Value = "ELC"
$ws.Range("G14").Value = "NRG"
$ws.Range("H14").Value = 0.8
$ws.Range("I14").NumberFormat = $DEC_FMT
$ws.Range("I14").Value = 29.89194313
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 1200
$ws.Range("M14").Value = 50
$ws.Range("O14").Value = 17

# --- Row 15: AuxStoOUT helper row (mirrors row 11) -------------------------
$ws.Range("F15").Value = "AuxStoOUT"
$ws.Range("G15").Formula = "=F15"
$ws.Range("I15").NumberFormat = $DEC_FMT
$ws.Range("L15").NumberFormat = $DEC_FMT
$ws.Range("L15").Formula = "=1/H14"

# --- Row 16: AuxStoIN helper row (mirrors row 12) --------------------------
$ws.Range("D16").NumberFormat = $TEXT_FMT
$ws.Range("F16").NumberFormat = $TEXT_FMT
$ws.Range("F16").Formula = "=C25"
$ws.Range("G16").Value = "ELC"
$ws.Range("I16").NumberFormat = $DEC_FMT
$ws.Range("P16").Value = 1

$ws.Range("J19").Select()
